$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the placeholder data row (templateCell1 / templateCell2 /
#    templateCell3) from the first table - only the header row remains.
# ---------------------------------------------------------------------
$t = $d.Tables(1)
if ($t.Rows.Count -ge 2) {
    $t.Rows(2).Delete()
}

# ---------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark that used to sit right after
#    "[dateApplication]" so that it now sits right after
#    "[totalFacturation]" instead (still collapsed, i.e. zero length).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$found = $d.Content.Duplicate
$found.Find.Execute("totalFacturation]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$found.Collapse(0)
$pos = $found.Start

# Zero-length ranges are not accepted reliably by Bookmarks.Add in this
# engine, so insert a throw-away marker character, wrap the bookmark
# around it, then delete the marker text - this leaves the
# bookmarkStart/bookmarkEnd pair immediately adjacent at the right spot.
$insPoint = $d.Range($pos, $pos)
$insPoint.InsertAfter("~")

$markerRange = $d.Range($pos, $pos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$markerRange2 = $d.Range($pos, $pos + 1)
$markerRange2.Text = ""
